# Apply the changes described in the commit diff:
#   1. Update the "as of" date in the confidential disclaimer text (cell A9)
#      from 2021-03-31 to 2021-04-05.
#   2. Update the Weight (column D) and Percent Change (column E) values
#      for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected, so it must be unprotected before any cell
# can be edited.
$ws.Unprotect()

# 1. Update the confidential disclaimer text in A9 (only the date changes).
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."

# Re-run autofit on row 9 so the row height used for the previous two-line
# text (which Excel recalculates whenever the text changes) goes back to
# the sheet's normal/automatic height instead of staying pinned to an
# explicit value.
$ws.Rows("9").AutoFit()

# 2. Update Weight / Percent Change values for rows 2-6.
$ws.Range("D2").Value = 0.2511067104489074
$ws.Range("E2").Value = 0.01232859479179771

$ws.Range("D3").Value = 0.2482762822281707
$ws.Range("E3").Value = 0.00783289817232391

$ws.Range("D4").Value = 0.2537600467525148
$ws.Range("E4").Value = 0.02066725715972839

$ws.Range("D5").Value = 0.2468569605704071
$ws.Range("E5").Value = 0.02105122016626448

$ws.Range("E6").Value = 0.01548168008956385

# Restore sheet protection with the same effective settings as before
# (contents & objects & scenarios protected, but column/row formatting
# still allowed).
$ws.Protect("lido", $true, $true, $true, $false, $false, $true, $true)
